# Merge "NGS tech & challenges" and "NGS Pipelines" sessions on 20/03/2024
# into a single longer "NGS tech & challenges" session (09:00-12:00), removing
# the separate "NGS Pipelines" slot. Everything below shifts up by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hyperlink objects in this engine don't re-anchor when rows shift, so drop
# them first and recreate them afterwards at their final locations.
$ws.Hyperlinks.Delete()

# "NGS tech & challenges" (row 15) now runs all the way to 12:00 (was 10:15),
# absorbing the slot previously used by "NGS Pipelines".
$ws.Range("D15").Value = 0.5

# Remove the old row 16 (NGS tech & challenges continuation, 10:15-11:00) and
# row 17 (NGS Pipelines, 11:00-12:00); rows below shift up by two.
$ws.Rows("16:17").Delete()

# Recreate the youtube hyperlinks at their final (shifted) positions.
$hyperlinkTargets = @{
    "J3"  = "https://youtu.be/qYp8shSMJ-0";
    "J6"  = "https://youtu.be/cxEtfKN91q4";
    "J8"  = "https://youtu.be/7MR1qUZQ94w";
    "J10" = "https://youtu.be/RGbONVWOaDo";
    "J12" = "https://youtu.be/gytTBNSWpFc";
    "J13" = "https://youtu.be/n3IpUHIodM8";
    "J15" = "https://youtu.be/5sUzkrucL1E";
    "J17" = "https://youtu.be/NxRECdxKP40";
    "J20" = "https://youtu.be/-cL0CI07-Es";
    "J23" = "https://youtu.be/gC_nslHgSa8";
    "J27" = "https://youtu.be/PfcrDlhY1zE";
    "J28" = "https://youtu.be/4HbSAEU5iBM";
}

$order = @("J3","J6","J8","J10","J12","J13","J15","J17","J20","J23","J27","J28")
foreach ($ref in $order) {
    $ws.Hyperlinks.Add($ws.Range($ref), $hyperlinkTargets[$ref])
}

# The named range "schedule" covered the whole table; shrink it to match the
# two fewer rows.
$wb.Names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$K`$28"

# Restore the author's last selected cell (shifted up along with the data).
$ws.Range("F31").Select()

Write-Output "Merged NGS sequencing and pipelines sessions"
